# Generate Report for Handback
#
# After a localization handback run completes, the status report workbook is
# refreshed for each locale:
#   - Status text moves from "Ready for handoff" to "Handed back: in sync
#     with en-US" (shown on the Overview sheet and on every locale sheet).
#   - The "Latest Target File" / "Latest Handback File" / "Latest Handback
#     DateTime" columns are populated for every source file, and a
#     hyperlink to the source markdown file is attached to the new
#     "Latest Target File" cell (matching the existing hyperlink on the
#     "Source File Name" cell).
#   - Columns that now hold longer text are widened so the new content is
#     readable.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$srcUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/72d1f373700320f253e7613395dc545f85a983c0/e2e/"

# -----------------------------------------------------------------------
# Overview sheet: status columns (zh-cn = E, de-de = F) for both rows.
# -----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $newStatus

# Widen the now-longer status columns to fit the new text.
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

# -----------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): update Status / Latest Target File /
# Latest Handback File / Latest Handback DateTime for each row, and add a
# hyperlink on the new Latest Target File cell.
# -----------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; HandbackDate = "2016-08-27 18:37:17" },
    @{ Sheet = "de-de"; HandbackDate = "2016-08-27 18:37:24" }
)

$rows = @(
    @{ Row = 2; SourceFile = "01eb3ba1-a6fc-4bc0-a925-cf4c0a11924c.md"; Hash = "9c134a94725f339cc56be235391b7b3c35638492" },
    @{ Row = 3; SourceFile = "66821751-ad8e-4a98-b6a4-5fda9c7ce2f3.md"; Hash = "31055bd7b67bea82edb99649a89fed46123ca439" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    foreach ($r in $rows) {
        $row = $r.Row
        $handbackFile = "$($r.SourceFile).$($r.Hash).$($locale.Sheet).xlf"

        # Status
        $ws.Cells.Item($row, 3).Value = $newStatus

        # Latest Target File (I) + hyperlink to the source file
        $targetCell = $ws.Cells.Item($row, 9)
        $targetCell.Value = $r.SourceFile
        $ws.Hyperlinks.Add($targetCell, "$srcUrlBase$($r.SourceFile)", $null, $null, $r.SourceFile) | Out-Null

        # Latest Handback File (J)
        $ws.Cells.Item($row, 10).Value = $handbackFile

        # Latest Handback DateTime (K)
        $ws.Cells.Item($row, 11).Value = $locale.HandbackDate
    }

    # Widen Status / Latest Target File / Latest Handback File columns.
    $ws.Columns.Item(3).ColumnWidth = 29.1666666666667
    $ws.Columns.Item(9).ColumnWidth = 39.1666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.1666666666667
}
